$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top of the data block (rows 2-8), shifting existing rows down
$ws.Rows("2:8").Insert()

# Clear any inherited formatting on the newly inserted rows so they match plain data rows
$ws.Rows("2:8").ClearFormats()

# Populate the 7 newly inserted rows with new sensor samples
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = 1.327703475952148
$ws.Cells.Item(2, 4).Value = -2.356002759933471
$ws.Cells.Item(2, 5).Value = -1.412894463539124
$ws.Cells.Item(2, 6).Value = 0.41813725233078
$ws.Cells.Item(2, 7).Value = -0.3617849349975586
$ws.Cells.Item(2, 8).Value = 0.7021896243095398

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = 1.294596195220947
$ws.Cells.Item(3, 4).Value = -2.42879855632782
$ws.Cells.Item(3, 5).Value = -1.538487493991852
$ws.Cells.Item(3, 6).Value = 0.2964223623275757
$ws.Cells.Item(3, 7).Value = 0.0287106670439243
$ws.Cells.Item(3, 8).Value = 0.1348485052585601

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = 1.5967435836792
$ws.Cells.Item(4, 4).Value = -2.721291160583497
$ws.Cells.Item(4, 5).Value = -1.309774732589722
$ws.Cells.Item(4, 6).Value = -0.015118914656341
$ws.Cells.Item(4, 7).Value = 0.0493273697793483
$ws.Cells.Item(4, 8).Value = 0.1067487001419067

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = 1.538975667953491
$ws.Cells.Item(5, 4).Value = -3.117227482795715
$ws.Cells.Item(5, 5).Value = -1.760656356811524
$ws.Cells.Item(5, 6).Value = -0.0705549344420433
$ws.Cells.Item(5, 7).Value = -0.0940732508897781
$ws.Cells.Item(5, 8).Value = -0.0430659987032413

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = 1.350284004211427
$ws.Cells.Item(6, 4).Value = -3.035107040405273
$ws.Cells.Item(6, 5).Value = -2.307204818725585
$ws.Cells.Item(6, 6).Value = 0.0039706239476799
$ws.Cells.Item(6, 7).Value = 0.1137736514210701
$ws.Cells.Item(6, 8).Value = 0.0995710343122482

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = 2.076921081542969
$ws.Cells.Item(7, 4).Value = -3.024871301651001
$ws.Cells.Item(7, 5).Value = -2.036388444900513
$ws.Cells.Item(7, 6).Value = 0.0445931628346443
$ws.Cells.Item(7, 7).Value = 0.0633772686123848
$ws.Cells.Item(7, 8).Value = 0.0119118718430399

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = 2.578349113464355
$ws.Cells.Item(8, 4).Value = -2.735702991485596
$ws.Cells.Item(8, 5).Value = -2.153444766998291
$ws.Cells.Item(8, 6).Value = -0.1004873365163803
$ws.Cells.Item(8, 7).Value = 0.0113010071218013
$ws.Cells.Item(8, 8).Value = 0.012980886735022

# Append 3 new rows of samples after the existing data (new rows 29-31)
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = 2.200880432128908
$ws.Cells.Item(29, 4).Value = -3.761867809295656
$ws.Cells.Item(29, 5).Value = -0.7659695267677282
$ws.Cells.Item(29, 6).Value = 0.1111774742603302
$ws.Cells.Item(29, 7).Value = -0.1542434692382812
$ws.Cells.Item(29, 8).Value = 0.3286454975605011

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = 2.088931465148925
$ws.Cells.Item(30, 4).Value = -3.697214221954345
$ws.Cells.Item(30, 5).Value = -0.6586695432662961
$ws.Cells.Item(30, 6).Value = 0.7035640478134155
$ws.Cells.Item(30, 7).Value = -0.6233879923820496
$ws.Cells.Item(30, 8).Value = 0.6962336897850037

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = 1.876247692108153
$ws.Cells.Item(31, 4).Value = -2.620123100280754
$ws.Cells.Item(31, 5).Value = -1.525603616237647
$ws.Cells.Item(31, 6).Value = 0.2559525370597839
$ws.Cells.Item(31, 7).Value = 0.1485929638147354
$ws.Cells.Item(31, 8).Value = -0.1579086631536483

# Re-stamp column A (elapsed-time index) for every data row, since the sliding window shift
# renumbers the time index sequentially regardless of which underlying sample moved where
$lastRow = $ws.UsedRange.Rows.Count
For ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}